# The deck ships two embedded DrawingML themes:
#   ppt/theme/theme1.xml -> "Office Theme" (default blue palette) - wired to the Notes Master
#   ppt/theme/theme2.xml -> "Integral"     (green palette)        - wired to the Slide Master
#       (and so to the presentation's visible design)
#
# The authored edit swaps which palette is which: the deck's visible design
# goes from "Integral" to the default "Office Theme" colors. The PowerPoint
# object model exposes the live theme as a single read/write 12-slot
# ThemeColorScheme (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink) on the
# slide/master, so recreate the swap by pushing the Office Theme palette's
# RGB values into each of those slots.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Office Theme palette (was theme1.xml) -- index order matches
# ThemeColorScheme.Colors(1..12): dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.
# RGB is passed as the VBA-style 0xBBGGRR packed integer.
$officeThemeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Colors($i).RGB = $officeThemeColors[$i - 1]
}
